# Daily attendance processing - 2026-01-24 11:56:23
# Normalize the "Recorded By" (column G) cell text: a few specific
# combinations of recorder names/emails need their comma-separated
# entries reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$map = @{
    "system, System, backup@backdoor.com" = "System, backup@backdoor.com, system"
    "dnasr281@gmail.com, System" = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
